# Auto-generated update script for cryptos worksheet
# Applies the per-cell price/volume changes described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.401.38'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').Value = '2.620.92'
$ws.Range('E3').Value = '  -2.02%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''111.15'
$ws.Range('E5').Value = '  -2.72%  '
$ws.Range('D6').Value = '''324.83'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('D7').Value = '''0.522'
$ws.Range('E7').Value = '  -1.72%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.541'
$ws.Range('E9').Value = '  -3.58%  '
$ws.Range('D10').Value = '''39.28'
$ws.Range('E10').Value = '  -5.16%  '
$ws.Range('D11').Value = '''19.90'
$ws.Range('E11').Value = '  -1.80%  '
$ws.Range('D12').Value = '''0.0806'
$ws.Range('E12').Value = '  -2.59%  '
$ws.Range('E13').Value = '  +1.21%  '
$ws.Range('D14').Value = '''7.30'
$ws.Range('E14').Value = '  -1.40%  '
$ws.Range('D15').Value = '3.032.73'
$ws.Range('E15').Value = '  -1.71%  '
$ws.Range('D16').Value = '2.617.12'
$ws.Range('E16').Value = '  -2.41%  '
$ws.Range('D17').Value = '''0.849'
$ws.Range('E17').Value = '  -3.60%  '
$ws.Range('D18').Value = '49.330.85'
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('E19').Value = '  -2.99%  '
$ws.Range('D20').Value = '''2.89'
$ws.Range('E20').Value = '  -1.63%  '
$ws.Range('D21').Value = '''6.65'
$ws.Range('E21').Value = '  -2.27%  '
$ws.Range('E22').Value = '  -2.64%  '
$ws.Range('D23').Value = '''265.41'
$ws.Range('E23').Value = '  -4.75%  '
$ws.Range('D24').Value = '''68.54'
$ws.Range('E24').Value = '  -5.92%  '
$ws.Range('E25').Value = '  -2.92%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '''25.88'
$ws.Range('E27').Value = '  -4.10%  '
$ws.Range('D28').Value = '''10.09'
$ws.Range('E28').Value = '  +1.19%  '
$ws.Range('D29').Value = '''2.19'
$ws.Range('E29').Value = '  -1.40%  '
$ws.Range('E30').Value = '  -3.11%  '
$ws.Range('D31').Value = '''34.33'
$ws.Range('E31').Value = '  -6.58%  '
$ws.Range('D32').Value = '''49.52'
$ws.Range('E32').Value = '  -2.01%  '
$ws.Range('D33').Value = '''5.47'
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('D34').Value = '''0.0805'
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').Value = '''18.78'
$ws.Range('E36').Value = '  -4.01%  '
$ws.Range('D37').Value = '''4.91'
$ws.Range('E37').Value = '  +1.94%  '
$ws.Range('D38').Value = '''2.02'
$ws.Range('E38').Value = '  -3.85%  '
$ws.Range('D39').Value = '''3.07'
$ws.Range('E39').Value = '  -1.02%  '
$ws.Range('D40').Value = '''128.14'
$ws.Range('E40').Value = '  +3.24%  '
$ws.Range('D41').Value = '''22.36'
$ws.Range('E41').Value = '  -1.44%  '
$ws.Range('D42').Value = '''0.111'
$ws.Range('E42').Value = '  -2.20%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '''2.20'
$ws.Range('E43').Value = '  -2.14%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '''0.0323'
$ws.Range('E44').Value = '  +2.23%  '
$ws.Range('D45').Value = '2.035.81'
$ws.Range('E45').Value = '  -2.31%  '
$ws.Range('E46').Value = '  +8.05%  '
$ws.Range('E47').Value = '  -5.65%  '
$ws.Range('E48').Value = '  -4.09%  '
$ws.Range('E49').Value = '  -3.86%  '
$ws.Range('D50').Value = '''5.17'
$ws.Range('E50').Value = '  -4.99%  '
$ws.Range('D51').Value = '''58.35'
$ws.Range('E51').Value = '  +0.90%  '
